$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.248.43'
$ws.Range('E2').Value = '  +3.69%  '
$ws.Range('D3').Value = '2.315.20'
$ws.Range('E3').Value = '  +2.32%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '518.76'
$ws.Range('E5').Value = '  +4.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.90'
$ws.Range('E6').Value = '  +4.16%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('E8').Value = '  +2.02%  '
$ws.Range('D9').Value = '2.337.39'
$ws.Range('E9').Value = '  +3.08%  '
$ws.Range('E10').Value = '  +8.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.155'
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.16'
$ws.Range('E12').Value = '  +7.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.343'
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '24.02'
$ws.Range('E14').Value = '  +4.71%  '
$ws.Range('D15').Value = '2.727.76'
$ws.Range('E15').Value = '  +2.36%  '
$ws.Range('D16').Value = '56.362.42'
$ws.Range('E16').Value = '  +3.89%  '
$ws.Range('E17').Value = '  +5.13%  '
$ws.Range('D18').Value = '2.333.98'
$ws.Range('E18').Value = '  +2.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.53'
$ws.Range('E19').Value = '  +3.02%  '
$ws.Range('E20').Value = '  +3.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '321.36'
$ws.Range('E21').Value = '  +6.27%  '
$ws.Range('E22').Value = '  +5.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.73'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('E25').Value = '  -0.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.159'
$ws.Range('E26').Value = '  +6.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.67'
$ws.Range('E27').Value = '  +5.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '171.67'
$ws.Range('E28').Value = '  +0.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.21'
$ws.Range('E29').Value = '  +11.93%  '
$ws.Range('E30').Value = '  +5.92%  '
$ws.Range('D31').Value = '0.0₃0731'
$ws.Range('E31').Value = '  +6.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.27'
$ws.Range('E32').Value = '  +5.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.35'
$ws.Range('E33').Value = '  +3.20%  '
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('E36').Value = '  +6.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.928'
$ws.Range('E37').Value = '  -0.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.02'
$ws.Range('E38').Value = '  +8.69%  '
$ws.Range('E39').Value = '  +9.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.50'
$ws.Range('E40').Value = '  +4.55%  '
$ws.Range('E41').Value = '  +2.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '140.85'
$ws.Range('E42').Value = '  +13.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.59'
$ws.Range('E43').Value = '  +6.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '280.65'
$ws.Range('E44').Value = '  +16.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.09'
$ws.Range('E45').Value = '  +6.07%  '
$ws.Range('E46').Value = '  +3.49%  '
$ws.Range('E47').Value = '  +3.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.558'
$ws.Range('E48').Value = '  +2.26%  '
$ws.Range('E49').Value = '  +2.40%  '
$ws.Range('E50').Value = '  +5.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.94'
$ws.Range('E51').Value = '  +5.36%  '
